$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The "transport@dacapo.com" aggregate-row (row 104) was recategorized
#    from department "transport" to "logistics" (the standalone "transport"
#    department value is retired).
$ws.Range("D104").Value = "logistics"

# 2. Append five new aggregated-email rows (106-110) to the lookup table
#    living in columns C:D, each with a mailto: hyperlink on the email cell
#    (matching the existing rows 102-105 pattern). Emails are filled in
#    first (column C, top to bottom), then the department column (D).
$newRows = @(
    @{ Row = 106; Email = "consnl@dacapo.com"; Dept = "warehouse " },
    @{ Row = 107; Email = "goodsreception@dacapo.com"; Dept = "warehouse" },
    @{ Row = 108; Email = "dkwm@dacapo.com"; Dept = "warehouse" },
    @{ Row = 109; Email = "ive@dacapo.com"; Dept = "warehouse" },
    @{ Row = 110; Email = "cutting@dacapo.com"; Dept = "production" }
)

foreach ($r in $newRows) {
    $cCell = $ws.Cells.Item($r.Row, 3)
    $cCell.Value = $r.Email
    $cCell.Style = "Hyperlink"
    $ws.Hyperlinks.Add($cCell, "mailto:" + $r.Email)
}

$ws.Cells.Item(106, 4).Value = "warehouse "
$ws.Cells.Item(107, 4).Value = "warehouse"
$ws.Cells.Item(108, 4).Value = "warehouse"
$ws.Cells.Item(109, 4).Value = "warehouse"

# 3. Row 74 (Gabriela Lasauri / Group Quality Manager) gains a note in the
#    new "additional information" column (E).
$ws.Range("E74").Value = "additional information"

$ws.Cells.Item(110, 4).Value = "production"

Write-Output "edit applied"
